$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 941.5714
$ws.Cells.Item(6, 9).Value = 222
$ws.Cells.Item(6, 11).Value = 666
$ws.Cells.Item(6, 13).Value = -554
$ws.Cells.Item(74, 8).Value = 9019.4
$ws.Cells.Item(74, 9).Value = 5999.5
$ws.Cells.Item(74, 10).Value = 9484
$ws.Cells.Item(74, 11).Value = 5999.5
$ws.Cells.Item(74, 12).Value = 9484
$ws.Cells.Item(74, 13).Value = -5063.5
$ws.Cells.Item(74, 14).Value = -11356
$ws.Cells.Item(77, 8).Value = 9019.4
$ws.Cells.Item(77, 9).Value = 5999.5
$ws.Cells.Item(77, 10).Value = 9484
$ws.Cells.Item(77, 11).Value = 29997.5
$ws.Cells.Item(77, 12).Value = 47420
$ws.Cells.Item(77, 13).Value = -25317.5
$ws.Cells.Item(77, 14).Value = -56780
$ws.Cells.Item(135, 8).Value = 1099.5
$ws.Cells.Item(135, 9).Value = 866.1667
$ws.Cells.Item(135, 10).Value = 1799.5
$ws.Cells.Item(135, 11).Value = 7795.5003
$ws.Cells.Item(135, 12).Value = 16195.5
$ws.Cells.Item(135, 13).Value = -5260.5003
$ws.Cells.Item(135, 14).Value = -21265.5
$ws.Cells.Item(137, 8).Value = 4097.294
$ws.Cells.Item(137, 9).Value = 1806.875
$ws.Cells.Item(137, 10).Value = 6133.222
$ws.Cells.Item(137, 11).Value = 5420.625
$ws.Cells.Item(137, 12).Value = 18399.666
$ws.Cells.Item(137, 13).Value = -2870.625
$ws.Cells.Item(137, 14).Value = -23499.666
$ws.Cells.Item(138, 8).Value = 4707.522
$ws.Cells.Item(138, 10).Value = 5665.6
$ws.Cells.Item(138, 12).Value = 16996.8
$ws.Cells.Item(138, 14).Value = -27276.8
$ws.Cells.Item(141, 8).Value = 2109.75
$ws.Cells.Item(141, 9).Value = 2246.3333
$ws.Cells.Item(141, 11).Value = 6738.999899999999
$ws.Cells.Item(141, 13).Value = -1558.999899999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 628.51135
$ws.Cells.Item(32, 9).Value = 566.77014
$ws.Cells.Item(32, 11).Value = 566.77014
$ws.Cells.Item(32, 13).Value = -279.77014
$ws.Cells.Item(63, 8).Value = 6107.3335
$ws.Cells.Item(63, 9).Value = 1596
$ws.Cells.Item(63, 10).Value = 7611.1113
$ws.Cells.Item(63, 11).Value = 1596
$ws.Cells.Item(63, 12).Value = 7611.1113
$ws.Cells.Item(63, 13).Value = -910
$ws.Cells.Item(63, 14).Value = -8983.1113
$ws.Cells.Item(66, 8).Value = 6107.3335
$ws.Cells.Item(66, 9).Value = 1596
$ws.Cells.Item(66, 10).Value = 7611.1113
$ws.Cells.Item(66, 11).Value = 7980
$ws.Cells.Item(66, 12).Value = 38055.5565
$ws.Cells.Item(66, 13).Value = -4548
$ws.Cells.Item(66, 14).Value = -44919.5565
$ws.Cells.Item(74, 8).Value = 7493.8887
$ws.Cells.Item(74, 9).Value = 4995
$ws.Cells.Item(74, 11).Value = 4995
$ws.Cells.Item(74, 13).Value = -4121
$ws.Cells.Item(77, 8).Value = 7493.8887
$ws.Cells.Item(77, 9).Value = 4995
$ws.Cells.Item(77, 11).Value = 24975
$ws.Cells.Item(77, 13).Value = -20607
$ws.Cells.Item(97, 8).Value = 964.6667
$ws.Cells.Item(97, 9).Value = 1043.8462
$ws.Cells.Item(97, 11).Value = 1043.8462
$ws.Cells.Item(97, 13).Value = -547.8462
$ws.Cells.Item(102, 8).Value = 64724.715
$ws.Cells.Item(102, 9).Value = 74695.664
$ws.Cells.Item(102, 11).Value = 74695.664
$ws.Cells.Item(102, 13).Value = -73073.664
$ws.Cells.Item(132, 8).Value = 2486.973
$ws.Cells.Item(132, 9).Value = 2493.6667
$ws.Cells.Item(132, 11).Value = 7481.000100000001
$ws.Cells.Item(132, 13).Value = -4951.000100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(48, 8).Value = 300000
$ws.Cells.Item(48, 10).Value = 300000
$ws.Cells.Item(48, 12).Value = 300000
$ws.Cells.Item(48, 14).Value = -300830
$ws.Cells.Item(64, 8).Value = 267.66666
$ws.Cells.Item(64, 10).Value = 276.5
$ws.Cells.Item(64, 12).Value = 276.5
$ws.Cells.Item(64, 14).Value = -726.5
$ws.Cells.Item(67, 8).Value = 267.66666
$ws.Cells.Item(67, 10).Value = 276.5
$ws.Cells.Item(67, 12).Value = 276.5
$ws.Cells.Item(67, 14).Value = -1836.5
$ws.Cells.Item(107, 8).Value = 2102.2354
$ws.Cells.Item(107, 9).Value = 2476.3157
$ws.Cells.Item(107, 11).Value = 2476.3157
$ws.Cells.Item(107, 13).Value = -556.3157000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 21953
$ws.Cells.Item(31, 9).Value = 1388.0714
$ws.Cells.Item(31, 11).Value = 1388.0714
$ws.Cells.Item(31, 13).Value = -1093.0714
$ws.Cells.Item(34, 8).Value = 21953
$ws.Cells.Item(34, 9).Value = 1388.0714
$ws.Cells.Item(34, 11).Value = 1388.0714
$ws.Cells.Item(34, 13).Value = -1186.0714
$ws.Cells.Item(99, 8).Value = 3975.3044
$ws.Cells.Item(99, 10).Value = 4151.6
$ws.Cells.Item(99, 12).Value = 4151.6
$ws.Cells.Item(99, 14).Value = -7147.6
$ws.Cells.Item(105, 8).Value = 31330.154
$ws.Cells.Item(105, 9).Value = 44758
$ws.Cells.Item(105, 10).Value = 1117.5
$ws.Cells.Item(105, 11).Value = 44758
$ws.Cells.Item(105, 12).Value = 1117.5
$ws.Cells.Item(105, 13).Value = -43011
$ws.Cells.Item(105, 14).Value = -4611.5
$ws.Cells.Item(126, 8).Value = 3975.3044
$ws.Cells.Item(126, 10).Value = 4151.6
$ws.Cells.Item(126, 12).Value = 12454.8
$ws.Cells.Item(126, 14).Value = -17394.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 1271.28
$ws.Cells.Item(129, 9).Value = 355.125
$ws.Cells.Item(129, 10).Value = 2900
$ws.Cells.Item(129, 11).Value = 1065.375
$ws.Cells.Item(129, 12).Value = 8700
$ws.Cells.Item(129, 13).Value = 3934.625
$ws.Cells.Item(129, 14).Value = -18700
$ws.Cells.Item(131, 8).Value = 2876.318
$ws.Cells.Item(131, 9).Value = 1239.375
$ws.Cells.Item(131, 11).Value = 3718.125
$ws.Cells.Item(131, 13).Value = 1321.875
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2943.75
$ws.Cells.Item(102, 9).Value = 2031
$ws.Cells.Item(102, 11).Value = 2031
$ws.Cells.Item(102, 13).Value = -409
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2853.9167
$ws.Cells.Item(7, 9).Value = 2842.5715
$ws.Cells.Item(7, 10).Value = 2869.8
$ws.Cells.Item(7, 11).Value = 2842.5715
$ws.Cells.Item(7, 12).Value = 2869.8
$ws.Cells.Item(7, 13).Value = -2730.5715
$ws.Cells.Item(7, 14).Value = -3093.8
$ws.Cells.Item(25, 8).Value = 5000
$ws.Cells.Item(25, 9).Value = 500
$ws.Cells.Item(25, 10).Value = 14000
$ws.Cells.Item(25, 11).Value = 500
$ws.Cells.Item(25, 12).Value = 14000
$ws.Cells.Item(25, 13).Value = -270
$ws.Cells.Item(25, 14).Value = -14460
$ws.Cells.Item(40, 8).Value = 2581.9614
$ws.Cells.Item(40, 9).Value = 2527.4783
$ws.Cells.Item(40, 11).Value = 2527.4783
$ws.Cells.Item(40, 13).Value = -2391.4783
$ws.Cells.Item(126, 8).Value = 2853.9167
$ws.Cells.Item(126, 9).Value = 2842.5715
$ws.Cells.Item(126, 10).Value = 2869.8
$ws.Cells.Item(126, 11).Value = 8527.7145
$ws.Cells.Item(126, 12).Value = 8609.400000000001
$ws.Cells.Item(126, 13).Value = -6057.7145
$ws.Cells.Item(126, 14).Value = -13549.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 19300.92
$ws.Cells.Item(126, 9).Value = 22244.428
$ws.Cells.Item(126, 10).Value = 3847.5
$ws.Cells.Item(126, 11).Value = 66733.284
$ws.Cells.Item(126, 12).Value = 11542.5
$ws.Cells.Item(126, 13).Value = -64263.284
$ws.Cells.Item(126, 14).Value = -16482.5
